# "Fix up students example"
#
# The AGE column value for Güvenç Attila (row 6) was mistakenly entered as
# the text "five" instead of the number 5. Correct it to a numeric value.
# Once no cell references the old "five" shared string it is dropped from
# the shared-strings table automatically, which is the effect captured by
# the diff to xl/sharedStrings.xml.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 5

# Reflect the author's last on-screen position/selection when they saved:
# scrolled down one row (A2 at the top) with E7 as the active/selected cell.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E7").Select()
